$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 7856.5713
$ws.Range("I19").Value = 6999.4
$ws.Range("K19").Value = 6999.4
$ws.Range("M19").Value = -6824.4

$ws.Range("H106").Value = 0
$ws.Range("I106").Value = 0
$ws.Range("K106").Value = 0
$ws.Range("M106").ClearContents()

$ws.Range("H107").Value = 12931917
$ws.Range("I107").Value = 5953416.5
$ws.Range("K107").Value = 5953416.5
$ws.Range("M107").Value = -5951496.5

$ws.Range("H128").Value = 89018
$ws.Range("J128").Value = 89018
$ws.Range("L128").Value = 89018
$ws.Range("N128").Value = -98978

$ws.Range("H129").Value = 1050.7084
$ws.Range("I129").Value = 652.94116
$ws.Range("K129").Value = 1958.82348
$ws.Range("M129").Value = 3041.17652

$ws.Range("H131").Value = 1719.3846
$ws.Range("I131").Value = 1022.4545
$ws.Range("K131").Value = 3067.3635
$ws.Range("M131").Value = 1972.6365

$ws.Range("H137").Value = 3833.2
$ws.Range("I137").Value = 3345.7693
$ws.Range("J137").Value = 7001.5
$ws.Range("K137").Value = 10037.3079
$ws.Range("L137").Value = 21004.5
$ws.Range("M137").Value = -7487.3079
$ws.Range("N137").Value = -26104.5

$ws.Range("H138").Value = 1474186
$ws.Range("I138").Value = 1662.5
$ws.Range("J138").Value = 2504952.5
$ws.Range("K138").Value = 4987.5
$ws.Range("L138").Value = 7514857.5
$ws.Range("M138").Value = 152.5
$ws.Range("N138").Value = -7525137.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4354547
$ws.Range("I32").Value = 4656038.5
$ws.Range("K32").Value = 4656038.5
$ws.Range("M32").Value = -4655751.5

$ws.Range("H44").Value = 67037
$ws.Range("J44").Value = 67037
$ws.Range("L44").Value = 67037
$ws.Range("N44").Value = -68013

$ws.Range("H122").Value = 3564.425
$ws.Range("I122").Value = 2243.4614
$ws.Range("K122").Value = 6730.3842
$ws.Range("M122").Value = -4280.3842

$ws.Range("H123").Value = 67464.5
$ws.Range("J123").Value = 67464.5
$ws.Range("L123").Value = 67464.5
$ws.Range("N123").Value = -77264.5

$ws.Range("H132").Value = 3909.7656
$ws.Range("I132").Value = 2836.8462
$ws.Range("K132").Value = 8510.5386
$ws.Range("M132").Value = -5980.5386

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H130").Value = 67811.42999999999
$ws.Range("J130").Value = 67811.42999999999
$ws.Range("L130").Value = 67811.42999999999
$ws.Range("N130").Value = -77851.42999999999

$ws.Range("H139").Value = 64999.668
$ws.Range("J139").Value = 67499.5
$ws.Range("L139").Value = 67499.5
$ws.Range("N139").Value = -77779.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6188.64
$ws.Range("I31").Value = 3544.5715
$ws.Range("K31").Value = 3544.5715
$ws.Range("M31").Value = -3249.5715

$ws.Range("H34").Value = 6188.64
$ws.Range("I34").Value = 3544.5715
$ws.Range("K34").Value = 3544.5715
$ws.Range("M34").Value = -3342.5715

$ws.Range("H51").Value = 33095.855
$ws.Range("J51").Value = 33095.855
$ws.Range("L51").Value = 33095.855
$ws.Range("N51").Value = -34567.855

$ws.Range("H58").Value = 4490.0967
$ws.Range("I58").Value = 1720.25
$ws.Range("K58").Value = 1720.25
$ws.Range("M58").Value = -1517.25

$ws.Range("H61").Value = 33095.855
$ws.Range("J61").Value = 33095.855
$ws.Range("L61").Value = 33095.855
$ws.Range("N61").Value = -33791.855

$ws.Range("H62").Value = 4778.7
$ws.Range("I62").Value = 4532.25
$ws.Range("J62").Value = 5148.375
$ws.Range("K62").Value = 4532.25
$ws.Range("L62").Value = 5148.375
$ws.Range("M62").Value = -3908.25
$ws.Range("N62").Value = -6396.375

$ws.Range("H65").Value = 4778.7
$ws.Range("I65").Value = 4532.25
$ws.Range("J65").Value = 5148.375
$ws.Range("K65").Value = 22661.25
$ws.Range("L65").Value = 25741.875
$ws.Range("M65").Value = -19541.25
$ws.Range("N65").Value = -31981.875

$ws.Range("H99").Value = 2744.6538
$ws.Range("I99").Value = 1486.7333
$ws.Range("K99").Value = 1486.7333
$ws.Range("M99").Value = 11.2666999999999

$ws.Range("H122").Value = 4159.5537
$ws.Range("I122").Value = 3494.6099
$ws.Range("K122").Value = 10483.8297
$ws.Range("M122").Value = -8033.8297

$ws.Range("H124").Value = 49868.5
$ws.Range("J124").Value = 49868.5
$ws.Range("L124").Value = 49868.5
$ws.Range("N124").Value = -54778.5

$ws.Range("H126").Value = 2744.6538
$ws.Range("I126").Value = 1486.7333
$ws.Range("K126").Value = 4460.199900000001
$ws.Range("M126").Value = -1990.199900000001

$ws.Range("H134").Value = 3752.9119
$ws.Range("I134").Value = 2184.16
$ws.Range("K134").Value = 6552.48
$ws.Range("M134").Value = -4017.48

$ws.Range("H136").Value = 4490.0967
$ws.Range("I136").Value = 1720.25
$ws.Range("K136").Value = 5160.75
$ws.Range("M136").Value = -2610.75

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 758348.75
$ws.Range("I4").Value = 2613.3333
$ws.Range("K4").Value = 7839.999899999999
$ws.Range("M4").Value = -7727.999899999999

$ws.Range("H63").Value = 12
$ws.Range("I63").Value = 12
$ws.Range("K63").Value = 36
$ws.Range("M63").Value = 713

$ws.Range("H66").Value = 12
$ws.Range("I66").Value = 12
$ws.Range("K66").Value = 108
$ws.Range("M66").Value = 3636

$ws.Range("H75").Value = 111121240
$ws.Range("I75").Value = 111111360
$ws.Range("J75").Value = 111131110
$ws.Range("K75").Value = 333334080
$ws.Range("L75").Value = 333393330
$ws.Range("M75").Value = -333333082
$ws.Range("N75").Value = -333395326

$ws.Range("H78").Value = 111121240
$ws.Range("I78").Value = 111111360
$ws.Range("J78").Value = 111131110
$ws.Range("K78").Value = 1000002240
$ws.Range("L78").Value = 1000179990
$ws.Range("M78").Value = -999997248
$ws.Range("N78").Value = -1000189974

$ws.Range("H131").Value = 2528.18
$ws.Range("I131").Value = 2580
$ws.Range("J131").Value = 2522.422
$ws.Range("K131").Value = 7740
$ws.Range("L131").Value = 7567.266
$ws.Range("M131").Value = -2700
$ws.Range("N131").Value = -17647.266

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 235.625
$ws.Range("I2").Value = 137.5
$ws.Range("K2").Value = 137.5
$ws.Range("M2").Value = -24.5

$ws.Range("H57").Value = 58294.5
$ws.Range("J57").Value = 66122.28999999999
$ws.Range("L57").Value = 66122.28999999999
$ws.Range("N57").Value = -67762.28999999999

$ws.Range("H102").Value = 3279.389
$ws.Range("I102").Value = 3106.121
$ws.Range("K102").Value = 3106.121
$ws.Range("M102").Value = -1484.121

$ws.Range("H122").Value = 1648682.2
$ws.Range("I122").Value = 2014338.5
$ws.Range("K122").Value = 6043015.5
$ws.Range("M122").Value = -6040565.5

$ws.Range("H126").Value = 19234260
$ws.Range("I126").Value = 45456596
$ws.Range("K126").Value = 136369788
$ws.Range("M126").Value = -136367318

$ws.Range("H132").Value = 2856.375
$ws.Range("I132").Value = 1635.4
$ws.Range("K132").Value = 4906.200000000001
$ws.Range("M132").Value = -2376.200000000001

$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5739.278
$ws.Range("I7").Value = 4649.6665
$ws.Range("J7").Value = 6284.0835
$ws.Range("K7").Value = 4649.6665
$ws.Range("L7").Value = 6284.0835
$ws.Range("M7").Value = -4537.6665
$ws.Range("N7").Value = -6508.0835

$ws.Range("H122").Value = 3737.3333
$ws.Range("I122").Value = 2733.6086
$ws.Range("K122").Value = 8200.825800000001
$ws.Range("M122").Value = -5750.825800000001

$ws.Range("H123").Value = 40429
$ws.Range("J123").Value = 40429
$ws.Range("L123").Value = 40429
$ws.Range("N123").Value = -50229

$ws.Range("H126").Value = 5739.278
$ws.Range("I126").Value = 4649.6665
$ws.Range("J126").Value = 6284.0835
$ws.Range("K126").Value = 13948.9995
$ws.Range("L126").Value = 18852.2505
$ws.Range("M126").Value = -11478.9995
$ws.Range("N126").Value = -23792.2505

$ws.Range("H136").Value = 7337.887
$ws.Range("I136").Value = 2148.4375
$ws.Range("K136").Value = 6445.3125
$ws.Range("M136").Value = -3895.3125

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 885.5454999999999
$ws.Range("I100").Value = 327.5
$ws.Range("J100").Value = 1204.4286
$ws.Range("K100").Value = 655
$ws.Range("L100").Value = 2408.8572
$ws.Range("M100").Value = -114
$ws.Range("N100").Value = -3490.8572

$ws.Range("H122").Value = 13625585
$ws.Range("I122").Value = 17381904
$ws.Range("K122").Value = 52145712
$ws.Range("M122").Value = -52143262

$ws.Range("H126").Value = 83337910
$ws.Range("I126").Value = 100002910
$ws.Range("K126").Value = 300008730
$ws.Range("M126").Value = -300006260

$ws.Range("H129").Value = 99000
$ws.Range("J129").Value = 99000
$ws.Range("L129").Value = 99000
$ws.Range("N129").Value = -109000

$ws.Range("H132").Value = 2859.0833
$ws.Range("I132").Value = 1288
$ws.Range("K132").Value = 3864
$ws.Range("M132").Value = -1334
